$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 104 - this shifts the existing rows 104-158
# down to 105-159 (and the sheet's used range grows from R158 to R159).
$ws.Rows.Item(104).Insert()

# Populate the newly inserted row 104 with the new weekly price record.
$ws.Range("A104").Value = 4
$ws.Range("B104").Value = 'Feria Lagunitas de Puerto Montt'
$ws.Range("C104").Value = 'Los Lagos'
$ws.Range("D104").Value = 44452
$ws.Range("E104").Value = 10
$ws.Range("F104").Value = 100112037
$ws.Range("G104").Value = 'Cebollín'
$ws.Range("H104").Value = 'Sin especificar'
$ws.Range("I104").Value = 'Primera'
$ws.Range("J104").Value = 80
$ws.Range("K104").Value = 6500
$ws.Range("L104").Value = 6500
$ws.Range("M104").Value = 6500
$ws.Range("N104").Value = '$/paquete 36 unidades'
$ws.Range("O104").Value = 'Región Metropolitana'
$ws.Range("P104").Value = 181
$ws.Range("Q104").Value = 36
$ws.Range("R104").Value = 'Hortaliza'
